$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLC Tags")

$ws.Range("E24").Value = "Transportband Werkstück am Bandanfang (100BG1)"
$ws.Range("E25").Value = "Transportband Werkstück am Vereinzeler (100BG2)"
$ws.Range("E26").Value = "Transportband Werkstück beim Bandende (100BG3)"
$ws.Range("E27").Value = "Pick&Placer Schlitten eingefahren (120BG1)"
$ws.Range("E28").Value = "Pick&Placer Schlitten ausgefahren (120BG2)"
$ws.Range("E29").Value = "Pick&Placer Sauger oben (120BG3)"
$ws.Range("E30").Value = "Pick&Placer Werkstück angesaugt (120PB4)"
$ws.Range("E31").Value = "Pick&Placer Bedienstelle Reset (160SF2)"
$ws.Range("E32").Value = "Pick&Placer Bedienstelle Start (160SF3)"
$ws.Range("E33").Value = "Pick&Placer Bedienstelle Stop (160SF4)"
$ws.Range("E34").Value = "Pick&Placer Bedienstelle Manual - Auto (160SF5)"
$ws.Range("E39").Value = "Transportband Vereinzeler ausfahren (105MB8)"
$ws.Range("E40").Value = "Pick&Placer Turmleuchte grün (130PF1)"
$ws.Range("E41").Value = "Pick&Placer Turmleuchte gelb (130PF1)"
$ws.Range("E42").Value = "Pick&Placer Turmleuchte rot (130PF1)"
$ws.Range("E43").Value = "Pick&Placer Schlitten einfahren (125MB1)"
$ws.Range("E44").Value = "Pick&Placer Schlitten ausfahren (125MB2)"
$ws.Range("E45").Value = "Pick&Placer Sauger nach unten (125MB3)"
$ws.Range("E46").Value = "Pick&Placer Vakuum ein (125MB4)"
$ws.Range("E47").Value = "Pick&Placer Bedienstelle LED Steuerspannung Ein (Reset) (160SF2)"
$ws.Range("E48").Value = "Pick&Placer Bedienstelle LED Fehler (160PF2)"
$ws.Range("E49").Value = "Pick&Placer Bedienstelle LED Start (160SF3)"
$ws.Range("E57").Value = "Transportband Bandmotor Geschwindigkeit Option (105TF1)"
$ws.Range("E59").Value = "Pick&Placer Bedienstelle Not-Halt Channel 1 (160SF1)"
$ws.Range("E63").Value = "Pick&Placer Bedienstelle Not-Halt Channel 2 (160SF1)"
$ws.Range("E67").Value = "Freigabe Not-Halt Channel 1 (81KF7)"
$ws.Range("E68").Value = "Freigabe Not-Halt Channel 2 (81KF8)"
$ws.Range("E69").Value = "Transportband Bandmotor Vorwärts (105TF1)"
$ws.Range("E70").Value = "Transportband Bandmotor Rückwärts (105TF1)"
